# Replace the four "Сазвежђе Персеус ..." paragraphs with the single
# translated-dates line "Сазвежђе Pegasus: 8. и 17. октобра, 7. и 16.
# новембра,", stripping all the old run formatting (the replacement text
# ends up in one bare, unformatted run) and dropping the leftover
# "_Hlk514861060" bookmark that used to wrap the first occurrence.

$d = $word.ActiveDocument

# The stray hidden bookmark around the first occurrence must go away
# entirely (both its start and end markers). Doing this before we touch
# the paragraph text keeps the engine's bookmark bookkeeping happy.
Try {
    $d.Bookmarks.Item("_Hlk514861060").Delete()
} Catch {
}

$newText = "Сазвежђе Pegasus: 8. и 17. октобра, 7. и 16. новембра,"

$searchRange = $d.Range(0, $d.Content.End)
while ($searchRange.Find.Execute("Сазвежђе", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)) {
    $para = $searchRange.Paragraphs(1)
    $paraRange = $para.Range

    # Exclude the trailing paragraph mark from the range we clear out.
    $target = $d.Range($paraRange.Start, $paraRange.End - 1)
    $target.Delete()
    $target.InsertAfter($newText)

    # Continue searching after the text we just inserted so we do not
    # re-match the "Сазвежђе" that is now part of the replacement text.
    $searchRange = $d.Range($target.End, $d.Content.End)
}
